$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Report (가로) sheet: add floatimage markers ---
$ws2.Range("A18").Value = '${floatimage(ci, "", 0:-1)}'
$ws2.Range("B30").Value = '${floatimage(ci, "", -1:-1)}'
$ws2.Range("E7").Value = '${floatimage(logo, "", 0:-0)}'
$ws2.Range("B20").Value = '${floatimage(logo, "", 0:-0)}'

$ws2.Range("B23").Value = '${floatimage(ci, "", 0:0)}'
[void]$ws2.Range("B23:C27").Merge()
$ws2.Range("B23:C27").HorizontalAlignment = -4108
$ws2.Range("B23:C27").VerticalAlignment = -4108

# --- Report(세로) sheet: add floatimage markers ---
$ws1.Range("B10").Value = '${floatimage(ci, B11, 0:-1)}'
$ws1.Range("E6").Value = '${floatimage(ci, C10, 0:-1)}'

# --- Update selections / active sheet to match latest editing session ---
[void]$ws2.Range("A18").Select()
[void]$ws3.Range("E7").Select()
$ws1.Activate()
[void]$ws1.Range("E7").Select()
